$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "PACS Submit Status" cell (O2) used to read "Successfully verify" (and the
# shared-string pool separately held stray "Verify" / "Successfully verify"
# entries). The new content for O2 is "Successfully Verified".
$ws.Range("O2").Value = "Successfully Verified"

# Column D ("District") widened from 12.0 to ~20.44 OOXML width units
# (the same visual width already used by column F).
$ws.Columns.Item(4).ColumnWidth = 19.6666666667

# The view was scrolled so the frozen/leftmost visible column moves from G to K.
$excel.ActiveWindow.ScrollColumn = 11
$excel.ActiveWindow.ScrollRow = 1
